$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D1 header text (identity/status field description)
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"

# Widen column D to fit the longer header text
$ws.Columns.Item(4).ColumnWidth = 26.375

# Update the active selection to F4
$ws.Range("F4").Select()
